# Automatische test-sync: 2025-08-28 18:38:50
# Append a new "Retour status" log entry to the Logs sheet and bump the
# Dashboard summary count for "Retour / Terugbetaling".

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 11

$logs.Cells.Item($newRow, 1).Value = "Retour status"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 4).Value = "Retour / Terugbetaling"
$logs.Cells.Item($newRow, 6).Value = "2025-08-28 18:38:08"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Nee"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

$dashboard.Range("B2").Value = 10

# Extend the conditional-formatting ranges (D/G/H/I/J columns) so they
# keep covering the newly appended data row.
$logs.Range("D2:D10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D11"))
$logs.Range("G2:G10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G11"))
$logs.Range("H2:H10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H11"))
$logs.Range("I2:I10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I11"))
$logs.Range("J2:J10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J11"))
